$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as text in the
# source data (e.g. "218.77"). Temporarily force the column to Text
# format so Excel does not silently convert these assignments to
# numbers, then strip the temporary format back off.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "89.800.10"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "3.217.34"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "218.04"
$ws.Range("E5").Value = "  +5.49%  "
$ws.Range("D6").Value = "630.54"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("D7").Value = "0.394"
$ws.Range("E7").Value = "  +7.25%  "
$ws.Range("D8").Value = "0.698"
$ws.Range("E8").Value = "  +6.28%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "3.215.13"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").Value = "0.574"
$ws.Range("E11").Value = "  +7.26%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "33.57"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "3.817.71"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").Value = "89.612.95"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "3.217.80"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "0.0000237"
$ws.Range("E19").Value = "  +83.36%  "
$ws.Range("D20").Value = "3.50"
$ws.Range("E20").Value = "  +18.81%  "
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").Value = "439.79"
$ws.Range("E22").Value = "  +6.63%  "
$ws.Range("D23").Value = "8.67"
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "5.27"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").Value = "12.03"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("D27").Value = "82.23"
$ws.Range("E27").Value = "  +12.14%  "
$ws.Range("D28").Value = "3.388.60"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "4.19"
$ws.Range("E32").Value = "  +39.28%  "
$ws.Range("D33").Value = "8.59"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").Value = "545.36"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "7.05"
$ws.Range("E35").Value = "  +6.74%  "
$ws.Range("E36").Value = "  +3.81%  "
$ws.Range("D37").Value = "1.32"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "0.130"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("D44").Value = "0.377"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "146.92"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "174.06"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "43.74"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.756"
$ws.Range("E48").Value = "  +10.08%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "1.26"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").Value = "0.125"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "0.626"
$ws.Range("E51").Value = "  +6.27%  "

$ws.Range("D2:D51").ClearFormats()
